# IS Team Updates - November 22
# Swap the display order of the "Picture only" slide and the "Today Schedule"
# speaker-list (table) slide, and refresh the speaker list table: insert a new
# speaker ("Connie Liu") after "John Hynes", pushing the rest of the agenda
# down a slot (and appending "Shazi" at the end), resize/reposition the table,
# and reduce its header/body font sizes.

$p = $ppt.ActivePresentation

# --- 1. Reorder slides -----------------------------------------------------
# The speaker-list table currently sits on slide 3, after the picture-only
# slide (slide 2). Move it to slide 2 so it now comes first.
$p.Slides.Item(3).MoveTo(2)

# --- 2. Locate the speaker-list table --------------------------------------
$s = $p.Slides.Item(2)
$tblShape = $null
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    if ($s.Shapes.Item($j).HasTable) {
        $tblShape = $s.Shapes.Item($j)
    }
}
$tbl = $tblShape.Table

# --- 3. Insert the new "Connie Liu" row after "John Hynes" (row 3) --------
$tbl.Rows.Add(4) | Out-Null
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Connie Liu"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "Presenting new team member"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "1 minute"

# --- 4. Shrink the header/body text ----------------------------------------
for ($c = 1; $c -le 3; $c++) {
    $tbl.Cell(1, $c).Shape.TextFrame.TextRange.Font.Size = 24
}
for ($r = 2; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $tbl.Cell($r, $c).Shape.TextFrame.TextRange.Font.Size = 20
    }
}

# --- 5. Reposition / resize the table and its columns/rows -----------------
$tblShape.Left = 24.406692913385825
$tblShape.Top = 113.78377952755906

$tbl.Columns.Item(1).Width = 201.43527559055119
$tbl.Columns.Item(2).Width = 557.2418897637796
$tbl.Columns.Item(3).Width = 163.7840157480315

$tbl.Rows.Item(1).Height = 39.02535433070866
for ($r = 2; $r -le ($tbl.Rows.Count - 1); $r++) {
    $tbl.Rows.Item($r).Height = 37.167007874015745
}
$tbl.Rows.Item($tbl.Rows.Count).Height = 39.02535433070866
